$d = $word.ActiveDocument

# The document currently ends with an empty paragraph right before the
# final section break. Add two new paragraphs after it:
#   1. "[PUMP:TBD:1]" - plain (Normal style) paragraph
#   2. "BOLUS:SRS:2"   - "List Bullet" style paragraph

$lastRange = $d.Paragraphs.Last.Range

$p1 = $d.Paragraphs.Add($lastRange)
$p1.Range.Text = "[PUMP:TBD:1]"

$p2 = $d.Paragraphs.Add($lastRange)
$p2.Range.Text = "BOLUS:SRS:2"
$p2.Style = "List Bullet"
